$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 13-26 reuse the same text values already present in rows 2-12,
# so we copy whole rows (preserving their "text" cell type/shared-string
# typing and default formatting) instead of re-typing the values, which
# would make Excel auto-convert the numeric-looking strings into numbers.
$ws.Range("A2:C2").Copy($ws.Range("A13:C13"))
$ws.Range("A2:C2").Copy($ws.Range("A14:C14"))
$ws.Range("A2:C2").Copy($ws.Range("A15:C15"))
$ws.Range("A2:C2").Copy($ws.Range("A16:C16"))
$ws.Range("A3:C3").Copy($ws.Range("A17:C17"))
$ws.Range("A4:C4").Copy($ws.Range("A18:C18"))
$ws.Range("A5:C5").Copy($ws.Range("A19:C19"))
$ws.Range("A6:C6").Copy($ws.Range("A20:C20"))
$ws.Range("A7:C7").Copy($ws.Range("A21:C21"))

# Row 22 matches row 8's Document/Operador but row 7's Registros value ("8").
$ws.Range("A8:B8").Copy($ws.Range("A22:B22"))
$ws.Range("C7").Copy($ws.Range("C22"))

$ws.Range("A9:C9").Copy($ws.Range("A23:C23"))
$ws.Range("A10:C10").Copy($ws.Range("A24:C24"))
$ws.Range("A11:C11").Copy($ws.Range("A25:C25"))
$ws.Range("A12:C12").Copy($ws.Range("A26:C26"))
